# Auto-generated: applies 2024-10-18 data update to violent-crime-ytd workbook
# Updates the '2024' (column K) figures across Citywide Totals, By Neighborhood,
# and each individual neighborhood sheet, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6179
$ws.Range('K3').Value = 6366
$ws.Range('K4').Value = 1334
$ws.Range('K5').Value = 452
$ws.Range('K6').Value = 7014
$ws.Range('K7').Value = 21345

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K3').Value = 72
$ws.Range('K7').Value = 271

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 427
$ws.Range('K4').Value = 76
$ws.Range('K6').Value = 472
$ws.Range('K7').Value = 1397

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 161
$ws.Range('K7').Value = 466

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 334
$ws.Range('K7').Value = 929

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K6').Value = 80
$ws.Range('K7').Value = 352

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 242
$ws.Range('K7').Value = 722

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 169
$ws.Range('K7').Value = 502

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 144
$ws.Range('K4').Value = 17
$ws.Range('K7').Value = 351

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 188
$ws.Range('K7').Value = 625
$ws.Range('K8').Value = 1397
$ws.Range('K10').Value = 122
$ws.Range('K11').Value = 397
$ws.Range('K15').Value = 218
$ws.Range('K19').Value = 624
$ws.Range('K20').Value = 509
$ws.Range('K21').Value = 69
$ws.Range('K23').Value = 218
$ws.Range('K24').Value = 63
$ws.Range('K25').Value = 103
$ws.Range('K27').Value = 204
$ws.Range('K29').Value = 1162
$ws.Range('K31').Value = 237
$ws.Range('K33').Value = 929
$ws.Range('K34').Value = 122
$ws.Range('K37').Value = 722
$ws.Range('K42').Value = 792
$ws.Range('K44').Value = 180
$ws.Range('K45').Value = 29
$ws.Range('K47').Value = 148
$ws.Range('K48').Value = 268
$ws.Range('K51').Value = 275
$ws.Range('K53').Value = 271
$ws.Range('K56').Value = 23
$ws.Range('K63').Value = 60
$ws.Range('K65').Value = 502
$ws.Range('K66').Value = 67
$ws.Range('K67').Value = 840
$ws.Range('K71').Value = 63
$ws.Range('K72').Value = 108
$ws.Range('K73').Value = 189
$ws.Range('K78').Value = 239
$ws.Range('K79').Value = 538
$ws.Range('K83').Value = 466
$ws.Range('K84').Value = 169
$ws.Range('K85').Value = 994
$ws.Range('K88').Value = 227
$ws.Range('K90').Value = 196
$ws.Range('K91').Value = 244
$ws.Range('K95').Value = 352
$ws.Range('K96').Value = 222
$ws.Range('K97').Value = 168
$ws.Range('K99').Value = 351
$ws.Range('K101').Value = 21345

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K6').Value = 81
$ws.Range('K7').Value = 237

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 232
$ws.Range('K3').Value = 304
$ws.Range('K7').Value = 840

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K2').Value = 58
$ws.Range('K7').Value = 169

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K4').Value = 57
$ws.Range('K6').Value = 329
$ws.Range('K7').Value = 1162

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K5').Value = 1
$ws.Range('K7').Value = 268

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 186
$ws.Range('K3').Value = 189
$ws.Range('K7').Value = 624

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K3').Value = 47
$ws.Range('K7').Value = 180

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K3').Value = 238
$ws.Range('K6').Value = 296
$ws.Range('K7').Value = 792

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K6').Value = 56
$ws.Range('K7').Value = 122

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 71
$ws.Range('K7').Value = 239

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K2').Value = 25
$ws.Range('K7').Value = 63

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K3').Value = 77
$ws.Range('K7').Value = 218

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 70
$ws.Range('K7').Value = 222

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K6').Value = 51
$ws.Range('K7').Value = 244

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('K6').Value = 41
$ws.Range('K7').Value = 69

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 179
$ws.Range('K4').Value = 33
$ws.Range('K5').Value = 17
$ws.Range('K7').Value = 538

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 165
$ws.Range('K7').Value = 509

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 204
$ws.Range('K3').Value = 207
$ws.Range('K6').Value = 169
$ws.Range('K7').Value = 625

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K3').Value = 33
$ws.Range('K7').Value = 122

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K3').Value = 36
$ws.Range('K7').Value = 103

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K6').Value = 49
$ws.Range('K7').Value = 148

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K3').Value = 55
$ws.Range('K7').Value = 218

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('K3').Value = 17
$ws.Range('K7').Value = 67

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K6').Value = 128
$ws.Range('K7').Value = 397

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K3').Value = 51
$ws.Range('K7').Value = 189

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K3').Value = 50
$ws.Range('K7').Value = 188

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K6').Value = 94
$ws.Range('K7').Value = 168

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K6').Value = 95
$ws.Range('K7').Value = 227

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K2').Value = 55
$ws.Range('K6').Value = 72
$ws.Range('K7').Value = 204

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K3').Value = 56
$ws.Range('K7').Value = 196

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 74
$ws.Range('K6').Value = 89
$ws.Range('K7').Value = 275

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 324
$ws.Range('K3').Value = 342
$ws.Range('K5').Value = 29
$ws.Range('K6').Value = 244
$ws.Range('K7').Value = 994

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K6').Value = 17
$ws.Range('K7').Value = 63

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K2').Value = 22
$ws.Range('K3').Value = 26
$ws.Range('K7').Value = 108

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('K2').Value = 7
$ws.Range('K7').Value = 29

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('K6').Value = 11
$ws.Range('K7').Value = 23

